$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: new ORM record for TUCMS.docx
$ws.Range("B9").Value = 22081808
$ws.Range("C9").Value = "TUCMS.docs"
$ws.Range("D9").Value = "Kaung Myat Bo"

# Copy date formatting (s="1", m/d/yyyy) from an existing dated row, then set the value
$ws.Range("E8").Copy($ws.Range("E9"))
$ws.Range("E9").Value = 43212

# Update the selection shown in the sheet view
$ws.Range("H13").Select()
